$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column D for rows 182-302: "No" -> "Yes"
$ws.Range("D182:D302").Value2 = "Yes"

# Move the selection/view: previously scrolled to A167 with C182 selected;
# now select E9 (which is within the default visible area, so Excel drops
# the explicit topLeftCell scroll-position attribute).
$ws.Range("E9").Select()
